$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value = 327.33334
$ws.Cells.Item(28, 10).Value = 136.5
$ws.Cells.Item(28, 12).Value = 136.5
$ws.Cells.Item(28, 14).Value = -1106.5
# Row 33
$ws.Cells.Item(33, 8).Value = 393.25
$ws.Cells.Item(33, 9).Value = 339.8095
$ws.Cells.Item(33, 11).Value = 339.8095
$ws.Cells.Item(33, 13).Value = -110.8095
# Row 58
$ws.Cells.Item(58, 8).Value = 643.6842
$ws.Cells.Item(58, 9).Value = 72.94118
$ws.Cells.Item(58, 10).Value = 5495
$ws.Cells.Item(58, 11).Value = 218.82354
$ws.Cells.Item(58, 12).Value = 16485
$ws.Cells.Item(58, 13).Value = -68.82354000000001
$ws.Cells.Item(58, 14).Value = -16785
# Row 61
$ws.Cells.Item(61, 8).Value = 1943.3334
$ws.Cells.Item(61, 9).Value = 415
$ws.Cells.Item(61, 11).Value = 1245
$ws.Cells.Item(61, 13).Value = -1073
# Row 81
$ws.Cells.Item(81, 8).Value = 32198
$ws.Cells.Item(81, 9).Value = 32198
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 32198
$ws.Cells.Item(81, 12).ClearContents()
$ws.Cells.Item(81, 14).Value = 0
$ws.Cells.Item(81, 13).Value = -31200
# Row 84
$ws.Cells.Item(84, 8).Value = 32198
$ws.Cells.Item(84, 9).Value = 32198
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 96594
$ws.Cells.Item(84, 12).ClearContents()
$ws.Cells.Item(84, 14).Value = 0
$ws.Cells.Item(84, 13).Value = -91602
# Row 98
$ws.Cells.Item(98, 8).Value = 294.4
$ws.Cells.Item(98, 9).Value = 294.4
$ws.Cells.Item(98, 11).Value = 294.4
$ws.Cells.Item(98, 13).Value = 1203.6
# Row 122
$ws.Cells.Item(122, 8).Value = 294.4
$ws.Cells.Item(122, 9).Value = 294.4
$ws.Cells.Item(122, 11).Value = 883.1999999999999
$ws.Cells.Item(122, 13).Value = 1566.8
# Row 132
$ws.Cells.Item(132, 8).Value = 39885.37
$ws.Cells.Item(132, 9).Value = 39885.37
$ws.Cells.Item(132, 11).Value = 119656.11
$ws.Cells.Item(132, 13).Value = -117126.11
# Row 138
$ws.Cells.Item(138, 8).Value = 2692.6296
$ws.Cells.Item(138, 9).Value = 796.55554
$ws.Cells.Item(138, 10).Value = 3640.6667
$ws.Cells.Item(138, 11).Value = 2389.66662
$ws.Cells.Item(138, 12).Value = 10922.0001
$ws.Cells.Item(138, 13).Value = 2750.33338
$ws.Cells.Item(138, 14).Value = -21202.0001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 2546.0352
$ws.Cells.Item(32, 9).Value = 1546.9111
$ws.Cells.Item(32, 11).Value = 1546.9111
$ws.Cells.Item(32, 13).Value = -1259.9111
# Row 45
$ws.Cells.Item(45, 8).Value = 2495.9546
$ws.Cells.Item(45, 9).Value = 1827.6428
$ws.Cells.Item(45, 10).Value = 3665.5
$ws.Cells.Item(45, 11).Value = 1827.6428
$ws.Cells.Item(45, 12).Value = 3665.5
$ws.Cells.Item(45, 13).Value = -1450.6428
$ws.Cells.Item(45, 14).Value = -4419.5
# Row 61
$ws.Cells.Item(61, 8).Value = 3370.7585
$ws.Cells.Item(61, 9).Value = 2791.95
$ws.Cells.Item(61, 11).Value = 2791.95
$ws.Cells.Item(61, 13).Value = -2579.95
# Row 132
$ws.Cells.Item(132, 8).Value = 17712.967
$ws.Cells.Item(132, 9).Value = 1211.2273
$ws.Cells.Item(132, 10).Value = 58050.555
$ws.Cells.Item(132, 11).Value = 3633.6819
$ws.Cells.Item(132, 12).Value = 174151.665
$ws.Cells.Item(132, 13).Value = -1103.6819
$ws.Cells.Item(132, 14).Value = -179211.665
# Row 136
$ws.Cells.Item(136, 8).Value = 3370.7585
$ws.Cells.Item(136, 9).Value = 2791.95
$ws.Cells.Item(136, 11).Value = 8375.849999999999
$ws.Cells.Item(136, 13).Value = -5825.849999999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 87
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 12).ClearContents()
$ws.Cells.Item(87, 14).Value = 0
# Row 90
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 12).ClearContents()
$ws.Cells.Item(90, 14).Value = 0
# Row 94
$ws.Cells.Item(94, 8).Value = 4420.467
$ws.Cells.Item(94, 9).Value = 1699.6
$ws.Cells.Item(94, 10).Value = 5780.9
$ws.Cells.Item(94, 11).Value = 1699.6
$ws.Cells.Item(94, 12).Value = 5780.9
$ws.Cells.Item(94, 13).Value = -1248.6
$ws.Cells.Item(94, 14).Value = -6682.9
# Row 107
$ws.Cells.Item(107, 8).Value = 1131.6154
$ws.Cells.Item(107, 9).Value = 971.1
$ws.Cells.Item(107, 11).Value = 971.1
$ws.Cells.Item(107, 13).Value = 948.9

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 8895.571
$ws.Cells.Item(31, 9).Value = 12782.8
$ws.Cells.Item(31, 10).Value = 3179.0588
$ws.Cells.Item(31, 11).Value = 12782.8
$ws.Cells.Item(31, 12).Value = 3179.0588
$ws.Cells.Item(31, 13).Value = -12487.8
$ws.Cells.Item(31, 14).Value = -3769.0588
# Row 34
$ws.Cells.Item(34, 8).Value = 8895.571
$ws.Cells.Item(34, 9).Value = 12782.8
$ws.Cells.Item(34, 10).Value = 3179.0588
$ws.Cells.Item(34, 11).Value = 12782.8
$ws.Cells.Item(34, 12).Value = 3179.0588
$ws.Cells.Item(34, 13).Value = -12580.8
$ws.Cells.Item(34, 14).Value = -3583.0588
# Row 132
$ws.Cells.Item(132, 8).Value = 23973.916
$ws.Cells.Item(132, 9).Value = 36738.855
$ws.Cells.Item(132, 11).Value = 110216.565
$ws.Cells.Item(132, 13).Value = -107686.565
# Row 134
$ws.Cells.Item(134, 8).Value = 825.5
$ws.Cells.Item(134, 9).Value = 781.875
$ws.Cells.Item(134, 11).Value = 2345.625
$ws.Cells.Item(134, 13).Value = 189.375

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Cells.Item(7, 8).Value = 67.5
$ws.Cells.Item(7, 9).Value = 35
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 105
$ws.Cells.Item(7, 12).Value = 300
$ws.Cells.Item(7, 13).Value = 7
$ws.Cells.Item(7, 14).Value = -524
# Row 19
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 13).ClearContents()
# Row 116
$ws.Cells.Item(116, 8).Value = 765.75
$ws.Cells.Item(116, 9).Value = 265.4
$ws.Cells.Item(116, 10).Value = 1599.6666
$ws.Cells.Item(116, 11).Value = 796.1999999999999
$ws.Cells.Item(116, 12).Value = 4798.9998
$ws.Cells.Item(116, 13).Value = 2645.8
$ws.Cells.Item(116, 14).Value = -11682.9998
# Row 131
$ws.Cells.Item(131, 8).Value = 823.5684
$ws.Cells.Item(131, 10).Value = 841.75824
$ws.Cells.Item(131, 12).Value = 2525.27472
$ws.Cells.Item(131, 14).Value = -12605.27472
# Row 141
$ws.Cells.Item(141, 8).Value = 2951.6
$ws.Cells.Item(141, 9).Value = 280
$ws.Cells.Item(141, 10).Value = 4732.6665
$ws.Cells.Item(141, 11).Value = 840
$ws.Cells.Item(141, 12).Value = 14197.9995
$ws.Cells.Item(141, 13).Value = 4340
$ws.Cells.Item(141, 14).Value = -24557.9995

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Cells.Item(97, 8).Value = 1653.8572
$ws.Cells.Item(97, 9).Value = 934.4
$ws.Cells.Item(97, 10).Value = 3452.5
$ws.Cells.Item(97, 11).Value = 934.4
$ws.Cells.Item(97, 12).Value = 3452.5
$ws.Cells.Item(97, 13).Value = -438.4
$ws.Cells.Item(97, 14).Value = -4444.5
# Row 102
$ws.Cells.Item(102, 8).Value = 2799.4075
$ws.Cells.Item(102, 9).Value = 2918.2856
$ws.Cells.Item(102, 10).Value = 2383.3333
$ws.Cells.Item(102, 11).Value = 2918.2856
$ws.Cells.Item(102, 12).Value = 2383.3333
$ws.Cells.Item(102, 13).Value = -1296.2856
$ws.Cells.Item(102, 14).Value = -5627.3333
# Row 122
$ws.Cells.Item(122, 8).Value = 2943.9583
$ws.Cells.Item(122, 10).Value = 4137.375
$ws.Cells.Item(122, 12).Value = 12412.125
$ws.Cells.Item(122, 14).Value = -17312.125
# Row 132
$ws.Cells.Item(132, 8).Value = 24875.727
$ws.Cells.Item(132, 9).Value = 1577.8667
$ws.Cells.Item(132, 11).Value = 4733.6001
$ws.Cells.Item(132, 13).Value = -2203.6001

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 123
$ws.Cells.Item(123, 8).Value = 30939.5
$ws.Cells.Item(123, 10).Value = 30939.5
$ws.Cells.Item(123, 12).Value = 30939.5
$ws.Cells.Item(123, 14).Value = -40739.5
# Row 132
$ws.Cells.Item(132, 8).Value = 671538.5600000001
$ws.Cells.Item(132, 9).Value = 1096835.5
$ws.Cells.Item(132, 11).Value = 3290506.5
$ws.Cells.Item(132, 13).Value = -3287976.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 114
$ws.Cells.Item(114, 8).Value = 34000
$ws.Cells.Item(114, 10).Value = 34000
$ws.Cells.Item(114, 12).Value = 34000
$ws.Cells.Item(114, 14).Value = -42678
# Row 132
$ws.Cells.Item(132, 8).Value = 1790.7812
$ws.Cells.Item(132, 9).Value = 1627.35
$ws.Cells.Item(132, 11).Value = 4882.049999999999
$ws.Cells.Item(132, 13).Value = -2352.049999999999
# Row 136
$ws.Cells.Item(136, 8).Value = 2483043.2
$ws.Cells.Item(136, 9).Value = 4608923.5
$ws.Cells.Item(136, 11).Value = 13826770.5
$ws.Cells.Item(136, 13).Value = -13824220.5
